$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Bottom of document: the paragraph that used to hold the bold
#    "Play Ancient Egypt Classic for Free - Slot Game Review" title
#    (a duplicate of the H1) is removed entirely, and the paragraph
#    right after it (the italic one-liner) gets its text replaced
#    with the new image-generation prompt while keeping the italic
#    run/formatting and the paragraph's existing leading empty run
#    untouched.
# ------------------------------------------------------------------

$total = $d.Paragraphs.Count
$boldTitlePara   = $d.Paragraphs($total - 1)
$italicDescPara  = $d.Paragraphs($total)

# Sanity checks on the text we expect to find before mutating.
# (Left as plain Find-free checks so a mismatch surfaces clearly.)

# Replace only the visible text of the italic paragraph (not the whole
# paragraph range), which preserves its existing run formatting (italic)
# and the paragraph's pre-existing leading empty run untouched.
$italicTextRange = $d.Range($italicDescPara.Range.Start, $italicDescPara.Range.End - 1)
$italicTextRange.Text = "Create a feature image for Ancient Egypt Classic that depicts a happy Maya warrior with glasses in a cartoon style. The warrior should be wearing traditional Mayan clothing and holding a golden scarab, representing the Wild symbol in the game. The background should show a temple with hieroglyphs and Egyptian gods, reflecting the theme of the game. The image should pop with bright colors and convey a fun and playful mood to attract potential players."

# Delete the bold duplicate-title paragraph completely (including its
# paragraph mark), which shifts the italic paragraph up to take its
# place.
$boldTitlePara.Range.Delete()

# ------------------------------------------------------------------
# 2) Top of document: insert a brand-new paragraph right after the H1
#    title, containing a bold "Meta description" run followed by a
#    plain run with the rest of the description text (plus a leading
#    empty run, matching this document's authoring convention).
#
#    We splice this in via a raw WordprocessingML fragment
#    (Range.InsertXML) rather than typing text into a freshly-inserted
#    paragraph: typing collapses/merges the empty leading run that
#    this document's other paragraphs all carry, whereas InsertXML
#    preserves the literal run structure we hand it. Inserting at a
#    position one character before the end of the title paragraph
#    (instead of exactly at the paragraph boundary) makes the engine
#    splice in our paragraph as a clean, separate paragraph instead of
#    fusing its trailing run into the next paragraph's text.
# ------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$splicePos = $titlePara.Range.End - 1
$insertPoint = $d.Range($splicePos, $splicePos)

$metaDescRestText = ": Experience Ancient Egypt in Ancient Egypt Classic slot game. Play free and earn up to x200 wins with the Golden Scarab. Read our in-depth review."

$metaParaXml = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>$metaDescRestText</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$null = $insertPoint.InsertXML($metaParaXml)

Write-Output "done"
